$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "gap" columns (gap_v = (#.Truck - BKS#.Trucks)/BKS#.Trucks) for each of
# --- the 6 time-limit blocks (30,25,20,15,10,5 minutes), inserted right after
# --- each block's "BKS #. Trucks" column. ---
$ws.Range("H3:H27").FormulaR1C1  = "=(RC[-4]-RC[-1])/RC[-1]"
$ws.Range("O3:O27").FormulaR1C1  = "=(RC[-4]-RC[-1])/RC[-1]"
$ws.Range("V3:V27").FormulaR1C1  = "=(RC[-4]-RC[-1])/RC[-1]"
$ws.Range("AC3:AC27").FormulaR1C1 = "=(RC[-4]-RC[-1])/RC[-1]"
$ws.Range("AJ3:AJ27").FormulaR1C1 = "=(RC[-4]-RC[-1])/RC[-1]"
$ws.Range("AQ3:AQ27").FormulaR1C1 = "=(RC[-4]-RC[-1])/RC[-1]"

# Row 28 already holds column averages for the other columns; extend it to the
# new gap columns too.
$ws.Range("H28").Formula  = "=AVERAGE(H3:H27)"
$ws.Range("O28").Formula  = "=AVERAGE(O3:O27)"
$ws.Range("V28").Formula  = "=AVERAGE(V3:V27)"
$ws.Range("AC28").Formula = "=AVERAGE(AC3:AC27)"
$ws.Range("AJ28").Formula = "=AVERAGE(AJ3:AJ27)"
$ws.Range("AQ28").Formula = "=AVERAGE(AQ3:AQ27)"

# --- New "fleetmin" summary table (rows 33-37, columns D-J) ---
# Shared-string insertion order matters (matches target: Cost, Vehicles,
# gap_c, gap_v, fleetmin).
$ws.Range("D34").Value = "Cost"
$ws.Range("D35").Value = "Vehicles"
$ws.Range("D36").Value = "gap_c"
$ws.Range("D37").Value = "gap_v"
$ws.Range("D33").Value = "fleetmin"

$ws.Range("E33").Value = 30
$ws.Range("F33").Value = 25
$ws.Range("G33").Value = 20
$ws.Range("H33").Value = 15
$ws.Range("I33").Value = 10
$ws.Range("J33").Value = 5

$ws.Range("E34").Formula = "=C29"
$ws.Range("F34").Formula = "=J29"
$ws.Range("G34").Formula = "=Q29"
$ws.Range("H34").Formula = "=X29"
$ws.Range("I34").Formula = "=AE29"
$ws.Range("J34").Formula = "=AL29"

$ws.Range("E35").Formula = "=D29"
$ws.Range("F35").Formula = "=K29"
$ws.Range("G35").Formula = "=R29"
$ws.Range("H35").Formula = "=Y29"
$ws.Range("I35").Formula = "=AF29"
$ws.Range("J35").Formula = "=AM29"

$ws.Range("E36").Value = 0.00904
$ws.Range("F36").Value = 0.01112
$ws.Range("G36").Value = 0.00876
$ws.Range("H36").Value = 0.01372
$ws.Range("I36").Value = 0.01225
$ws.Range("J36").Value = 0.0111

$ws.Range("E37").Value = 0.10082
$ws.Range("F37").Value = 0.12108
$ws.Range("G37").Value = 0.11936
$ws.Range("H37").Value = 0.12508
$ws.Range("I37").Value = 0.12774
$ws.Range("J37").Value = 0.10149

$ws.Range("E36:J37").Style = "Percent"

$ws.Columns("I:J").ColumnWidth = 9.2857142857

$ws.Range("K43").Select() | Out-Null
